# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 209-210 (pushing the existing
# rows 209..298 down to 211..300), then populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 209 - this shifts the
# existing data (rows 209..298) down to rows 211..300.
$ws.Range("A209:T210").Insert()

# New row 209: Lane Late / Primera
$ws.Cells.Item(209, 1).Value = 4
$ws.Cells.Item(209, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(209, 3).Value = 'Los Lagos'
$ws.Cells.Item(209, 4).Value = 44523
$ws.Cells.Item(209, 5).Value = 10
$ws.Cells.Item(209, 6).Value = 'Fruta'
$ws.Cells.Item(209, 7).Value = 100102
$ws.Cells.Item(209, 8).Value = 'Cítricos'
$ws.Cells.Item(209, 9).Value = 100102005
$ws.Cells.Item(209, 10).Value = 'Naranja'
$ws.Cells.Item(209, 11).Value = 'Lane Late'
$ws.Cells.Item(209, 12).Value = 'Primera'
$ws.Cells.Item(209, 13).Value = 600
$ws.Cells.Item(209, 14).Value = 14000
$ws.Cells.Item(209, 15).Value = 14500
$ws.Cells.Item(209, 16).Value = 14250
$ws.Cells.Item(209, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(209, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(209, 19).Value = 950
$ws.Cells.Item(209, 20).Value = 15

# New row 210: Navel Late / Primera
$ws.Cells.Item(210, 1).Value = 4
$ws.Cells.Item(210, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(210, 3).Value = 'Los Lagos'
$ws.Cells.Item(210, 4).Value = 44523
$ws.Cells.Item(210, 5).Value = 10
$ws.Cells.Item(210, 6).Value = 'Fruta'
$ws.Cells.Item(210, 7).Value = 100102
$ws.Cells.Item(210, 8).Value = 'Cítricos'
$ws.Cells.Item(210, 9).Value = 100102005
$ws.Cells.Item(210, 10).Value = 'Naranja'
$ws.Cells.Item(210, 11).Value = 'Navel Late'
$ws.Cells.Item(210, 12).Value = 'Primera'
$ws.Cells.Item(210, 13).Value = 600
$ws.Cells.Item(210, 14).Value = 14000
$ws.Cells.Item(210, 15).Value = 14500
$ws.Cells.Item(210, 16).Value = 14250
$ws.Cells.Item(210, 17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(210, 18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(210, 19).Value = 950
$ws.Cells.Item(210, 20).Value = 15
